$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.967257559299469
$ws.Range("B1").Value = 1.186035394668579
$ws.Range("C1").Value = 1.008044481277466
$ws.Range("D1").Value = 1.041595458984375
$ws.Range("E1").Value = 1.185782074928284
